# LV_AC_Isolator_19DX101 "QL" sheet update.
#
# The 50A rating group (rows 19-22: 1P/2P/3P/4P, width=50) is removed from
# the table. Everything below it (63A, 80A, 100A, ... 6300A groups) shifts
# up by four rows, which is exactly what a native row delete does - all
# the now-updated values (D19=63, D23=80, D27=100/E27="22", etc. down to
# the last surviving row, the first 6300A group at row 63/64) fall out of
# this single operation with no per-cell edits required.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("19:22").Delete()

# Leave the view parked where the author last left it.
$ws.Range("N25").Select()
